$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 22153.285
$ws.Range("I6").Value = 22153.285
$ws.Range("K6").Value = 66459.855
$ws.Range("M6").Value = -66347.855
# Row 11
$ws.Range("H11").Value = 137.44444
$ws.Range("I11").Value = 137.44444
$ws.Range("K11").Value = 137.44444
$ws.Range("M11").Value = 2.555560000000014
# Row 19
$ws.Range("H19").Value = 1446.875
$ws.Range("I19").Value = 1449.75
$ws.Range("K19").Value = 1449.75
$ws.Range("M19").Value = -1274.75
# Row 70
$ws.Range("H70").Value = 3609.261
$ws.Range("I70").Value = 1272.3334
$ws.Range("K70").Value = 3817.0002
$ws.Range("M70").Value = -3547.0002
# Row 73
$ws.Range("H73").Value = 3609.261
$ws.Range("I73").Value = 1272.3334
$ws.Range("K73").Value = 3817.0002
$ws.Range("M73").Value = -2881.0002
# Row 116
$ws.Range("H116").Value = 6703.3687
$ws.Range("I116").Value = 6676.4688
$ws.Range("J116").Value = 6846.8335
$ws.Range("K116").Value = 6676.4688
$ws.Range("L116").Value = 6846.8335
$ws.Range("M116").Value = -3234.4688
$ws.Range("N116").Value = -13730.8335
# Row 118
$ws.Range("H118").Value = 1124
$ws.Range("I118").Value = 1124
$ws.Range("K118").Value = 3372
$ws.Range("M118").Value = -1715
# Row 125
$ws.Range("H125").Value = 8110.0625
$ws.Range("J125").Value = 9708.25
$ws.Range("L125").Value = 87374.25
$ws.Range("N125").Value = -92294.25
# Row 132
$ws.Range("H132").Value = 2754.5625
$ws.Range("I132").Value = 2888.9644
$ws.Range("J132").Value = 1813.75
$ws.Range("K132").Value = 8666.893199999999
$ws.Range("L132").Value = 5441.25
$ws.Range("M132").Value = -6136.893199999999
$ws.Range("N132").Value = -10501.25
# Row 137
$ws.Range("H137").Value = 2543.5405
$ws.Range("I137").Value = 2348.7334
$ws.Range("K137").Value = 7046.2002
$ws.Range("M137").Value = -4496.2002
# Row 138
$ws.Range("H138").Value = 2217.868
$ws.Range("I138").Value = 1431.25
$ws.Range("J138").Value = 2558.027
$ws.Range("K138").Value = 4293.75
$ws.Range("L138").Value = 7674.081
$ws.Range("M138").Value = 846.25
$ws.Range("N138").Value = -17954.081
# Row 141
$ws.Range("H141").Value = 2778.25
$ws.Range("I141").Value = 2778.25
$ws.Range("K141").Value = 8334.75
$ws.Range("M141").Value = -3154.75

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3106.3215
$ws.Range("I32").Value = 1518.0962
$ws.Range("K32").Value = 1518.0962
$ws.Range("M32").Value = -1231.0962
# Row 45
$ws.Range("H45").Value = 5701.5625
$ws.Range("I45").Value = 7046.4546
$ws.Range("J45").Value = 2742.8
$ws.Range("K45").Value = 7046.4546
$ws.Range("L45").Value = 2742.8
$ws.Range("M45").Value = -6669.4546
$ws.Range("N45").Value = -3496.8
# Row 61
$ws.Range("H61").Value = 250001500
$ws.Range("I61").Value = 250001500
$ws.Range("K61").Value = 250001500
$ws.Range("M61").Value = -250001288
# Row 74
$ws.Range("H74").Value = 41672732
$ws.Range("I74").Value = 50006456
$ws.Range("J74").Value = 4124
$ws.Range("K74").Value = 50006456
$ws.Range("L74").Value = 4124
$ws.Range("M74").Value = -50005582
$ws.Range("N74").Value = -5872
# Row 77
$ws.Range("H77").Value = 41672732
$ws.Range("I77").Value = 50006456
$ws.Range("J77").Value = 4124
$ws.Range("K77").Value = 250032280
$ws.Range("L77").Value = 20620
$ws.Range("M77").Value = -250027912
$ws.Range("N77").Value = -29356
# Row 102
$ws.Range("H102").Value = 4546980
$ws.Range("I102").Value = 7144102.5
$ws.Range("K102").Value = 7144102.5
$ws.Range("M102").Value = -7142480.5
# Row 136
$ws.Range("H136").Value = 250001500
$ws.Range("I136").Value = 250001500
$ws.Range("K136").Value = 750004500
$ws.Range("M136").Value = -750001950

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2912.3333
$ws.Range("I94").Value = 2912.3333
$ws.Range("K94").Value = 2912.3333
$ws.Range("M94").Value = -2461.3333
# Row 140
$ws.Range("H140").Value = 98697.5
$ws.Range("J140").Value = 98697.5
$ws.Range("L140").Value = 98697.5
$ws.Range("N140").Value = -109057.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 132
$ws.Range("H132").Value = 45456830
$ws.Range("I132").Value = 52633836
$ws.Range("J132").Value = 2466
$ws.Range("K132").Value = 157901508
$ws.Range("L132").Value = 7398
$ws.Range("M132").Value = -157898978
$ws.Range("N132").Value = -12458

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5942.2144
$ws.Range("I3").Value = 4432.6665
$ws.Range("K3").Value = 13297.9995
$ws.Range("M3").Value = -13185.9995
# Row 32
$ws.Range("H32").Value = 333998.16
$ws.Range("J32").Value = 500499
$ws.Range("L32").Value = 1501497
$ws.Range("N32").Value = -1502063
# Row 34
$ws.Range("H34").Value = 1174
$ws.Range("I34").Value = 965.6667
$ws.Range("K34").Value = 2897.0001
$ws.Range("M34").Value = -2813.0001
# Row 39
$ws.Range("H39").Value = 958.44446
$ws.Range("I39").Value = 704.5
$ws.Range("K39").Value = 2113.5
$ws.Range("M39").Value = -1819.5
# Row 55
$ws.Range("H55").Value = 999
$ws.Range("I55").Value = 999
$ws.Range("K55").Value = 2997
$ws.Range("M55").Value = -2820
# Row 69
$ws.Range("H69").Value = 494.5
$ws.Range("I69").Value = 494.5
$ws.Range("K69").Value = 1483.5
$ws.Range("M69").Value = -672.5
# Row 72
$ws.Range("H72").Value = 494.5
$ws.Range("I72").Value = 494.5
$ws.Range("K72").Value = 4450.5
$ws.Range("M72").Value = -394.5
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 129
$ws.Range("H129").Value = 3294.8
$ws.Range("J129").Value = 3626.9333
$ws.Range("L129").Value = 10880.7999
$ws.Range("N129").Value = -20880.7999

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 4000427.8
$ws.Range("I3").Value = 4000427.8
$ws.Range("K3").Value = 4000427.8
$ws.Range("M3").Value = -4000311.8
# Row 22
$ws.Range("H22").Value = 51000
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1471
# Row 24
$ws.Range("H24").Value = 19100
$ws.Range("J24").Value = 100
$ws.Range("L24").Value = 100
$ws.Range("N24").Value = -446
# Row 80
$ws.Range("H80").Value = 2298.625
$ws.Range("I80").Value = 2270.5715
$ws.Range("J80").Value = 2495
$ws.Range("K80").Value = 2270.5715
$ws.Range("L80").Value = 2495
$ws.Range("M80").Value = -1272.5715
$ws.Range("N80").Value = -4491
# Row 83
$ws.Range("H83").Value = 2298.625
$ws.Range("I83").Value = 2270.5715
$ws.Range("J83").Value = 2495
$ws.Range("K83").Value = 11352.8575
$ws.Range("L83").Value = 2495
$ws.Range("M83").Value = -6360.8575
$ws.Range("N83").Value = -22459
# Row 126
$ws.Range("H126").Value = 4098.75
$ws.Range("I126").Value = 4098.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12296.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9826.25
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 5438228.5
$ws.Range("I132").Value = 5438228.5
$ws.Range("K132").Value = 16314685.5
$ws.Range("M132").Value = -16312155.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2711.6428
$ws.Range("I7").Value = 2214.7273
$ws.Range("J7").Value = 4533.6665
$ws.Range("K7").Value = 2214.7273
$ws.Range("L7").Value = 4533.6665
$ws.Range("M7").Value = -2102.7273
$ws.Range("N7").Value = -4757.6665
# Row 126
$ws.Range("H126").Value = 2711.6428
$ws.Range("I126").Value = 2214.7273
$ws.Range("J126").Value = 4533.6665
$ws.Range("K126").Value = 6644.1819
$ws.Range("L126").Value = 13600.9995
$ws.Range("M126").Value = -4174.1819
$ws.Range("N126").Value = -18540.9995
# Row 136
$ws.Range("H136").Value = 2164.0417
$ws.Range("I136").Value = 1686.25
$ws.Range("K136").Value = 5058.75
$ws.Range("M136").Value = -2508.75

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 27000
$ws.Range("J82").Value = 27000
$ws.Range("L82").Value = 27000
$ws.Range("N82").Value = -27766
# Row 85
$ws.Range("H85").Value = 27000
$ws.Range("J85").Value = 27000
$ws.Range("L85").Value = 27000
